# Automatische test-sync: 2025-08-05 17:58:50
# Append a new log row (row 20) to the "Logs" sheet, extend the conditional
# formatting ranges that covered rows 2:19 to cover rows 2:20, and bump the
# "Planning / Afspraak" tally on the "Dashboard" sheet from 13 to 14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$newRow = 20

$ws.Cells.Item($newRow, 1).Value  = "Kun jij dit even regelen?"
$ws.Cells.Item($newRow, 2).Value  = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 3).Value  = "Testmail #1: Kun jij dit even regelen?"
$ws.Cells.Item($newRow, 4).Value  = "Planning / Afspraak"
$ws.Cells.Item($newRow, 5).Value  = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$ws.Cells.Item($newRow, 6).Value  = "2025-08-05 17:57:59"
$ws.Cells.Item($newRow, 7).Value  = "Ja"
$ws.Cells.Item($newRow, 8).Value  = "Ja"
$ws.Cells.Item($newRow, 9).Value  = "Nee"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# The sheet's used range (and <dimension>) now runs through row 20 instead
# of row 19. Extend every conditional-formatting rule that applied to
# rows 2:19 of columns D, G, H, I, J so it also covers the new row, while
# preserving rule order/priority and referenced dxf styles.
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $ws.Range($col + "2:" + $col + "19")
    $newRange = $ws.Range($col + "2:" + $col + "20")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Dashboard: the "Planning / Afspraak" count goes from 13 to 14 now that a
# new row with that category was logged.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(2, 2).Value = 14
